$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1111111111111111
$ws.Range("C2").Value = 0.7171717171717171
$ws.Range("P2").Value = 0.08080808080808081
$ws.Range("S2").Value = 0.09090909090909091
$ws.Range("B3").Value = 0.006896551724137931
$ws.Range("C3").Value = 0.01379310344827586
$ws.Range("J3").Value = 0.01379310344827586
$ws.Range("P3").Value = 0.7724137931034483
$ws.Range("S3").Value = 0.1931034482758621
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6170212765957447
$ws.Range("S4").Value = 0.3617021276595745
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.03017241379310345
$ws.Range("D6").Value = 0.004310344827586207
$ws.Range("F6").Value = 0.04741379310344827
$ws.Range("J6").Value = 0.228448275862069
$ws.Range("O6").Value = 0.01724137931034483
$ws.Range("Q6").Value = 0.1293103448275862
$ws.Range("R6").Value = 0.09913793103448276
$ws.Range("S6").Value = 0.4439655172413793
$ws.Range("B7").Value = 0.07692307692307693
$ws.Range("D7").Value = 0.04142011834319527
$ws.Range("F7").Value = 0.0650887573964497
$ws.Range("J7").Value = 0.02958579881656805
$ws.Range("O7").Value = 0.005917159763313609
$ws.Range("Q7").Value = 0.1597633136094675
$ws.Range("R7").Value = 0.09467455621301775
$ws.Range("S7").Value = 0.5266272189349113
$ws.Range("B8").Value = 0.06329113924050633
$ws.Range("D8").Value = 0.02531645569620253
$ws.Range("F8").Value = 0.06835443037974684
$ws.Range("J8").Value = 0.04050632911392405
$ws.Range("O8").Value = 0.01518987341772152
$ws.Range("Q8").Value = 0.1746835443037975
$ws.Range("R8").Value = 0.1139240506329114
$ws.Range("S8").Value = 0.4987341772151899
$ws.Range("B9").Value = 0.06363636363636363
$ws.Range("D9").Value = 0.01363636363636364
$ws.Range("F9").Value = 0.03181818181818181
$ws.Range("J9").Value = 0.07272727272727272
$ws.Range("O9").Value = 0.02727272727272727
$ws.Range("Q9").Value = 0.1863636363636364
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.5045454545454545
$ws.Range("B10").Value = 0.1062176165803109
$ws.Range("D10").Value = 0.02072538860103627
$ws.Range("E10").Value = 0.002590673575129534
$ws.Range("F10").Value = 0.06088082901554404
$ws.Range("J10").Value = 0.07772020725388601
$ws.Range("O10").Value = 0.009067357512953367
$ws.Range("Q10").Value = 0.2383419689119171
$ws.Range("R10").Value = 0.08549222797927461
$ws.Range("S10").Value = 0.3989637305699482
$ws.Range("G11").Value = 0.1439688715953307
$ws.Range("J11").Value = 0.0311284046692607
$ws.Range("K11").Value = 0.1867704280155642
$ws.Range("L11").Value = 0.5447470817120622
$ws.Range("S11").Value = 0.0933852140077821
$ws.Range("G12").Value = 0.7588652482269503
$ws.Range("J12").Value = 0.0851063829787234
$ws.Range("K12").Value = 0.02836879432624113
$ws.Range("L12").Value = 0.02127659574468085
$ws.Range("S12").Value = 0.1063829787234043
$ws.Range("F15").Value = 0.04504504504504504
$ws.Range("H15").Value = 0.1621621621621622
$ws.Range("I15").Value = 0.06756756756756757
$ws.Range("J15").Value = 0.2702702702702703
$ws.Range("K15").Value = 0.04504504504504504
$ws.Range("M15").Value = 0.01801801801801802
$ws.Range("O15").Value = 0.04954954954954955
$ws.Range("S15").Value = 0.3423423423423423
$ws.Range("F16").Value = 0.02994011976047904
$ws.Range("H16").Value = 0.1317365269461078
$ws.Range("I16").Value = 0.0718562874251497
$ws.Range("J16").Value = 0.2634730538922156
$ws.Range("K16").Value = 0.1137724550898204
$ws.Range("M16").Value = 0.01796407185628742
$ws.Range("O16").Value = 0.0718562874251497
$ws.Range("S16").Value = 0.2994011976047904
$ws.Range("F17").Value = 0.01720430107526882
$ws.Range("H17").Value = 0.1483870967741935
$ws.Range("I17").Value = 0.08817204301075268
$ws.Range("J17").Value = 0.2451612903225806
$ws.Range("K17").Value = 0.09032258064516129
$ws.Range("M17").Value = 0.02150537634408602
$ws.Range("O17").Value = 0.07956989247311828
$ws.Range("S17").Value = 0.3096774193548387
$ws.Range("F18").Value = 0.02912621359223301
$ws.Range("H18").Value = 0.1796116504854369
$ws.Range("I18").Value = 0.0825242718446602
$ws.Range("J18").Value = 0.2961165048543689
$ws.Range("K18").Value = 0.07766990291262135
$ws.Range("M18").Value = 0.009708737864077669
$ws.Range("O18").Value = 0.07281553398058252
$ws.Range("S18").Value = 0.2524271844660194
$ws.Range("F19").Value = 0.03815028901734104
$ws.Range("H19").Value = 0.1364161849710983
$ws.Range("I19").Value = 0.07745664739884393
$ws.Range("J19").Value = 0.1890173410404624
$ws.Range("K19").Value = 0.06358381502890173
$ws.Range("M19").Value = 0.01271676300578035
$ws.Range("N19").Value = 0.0005780346820809249
$ws.Range("O19").Value = 0.05953757225433526
$ws.Range("S19").Value = 0.4225433526011561
